$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the old row 16 (the blank separator row right
# before the USB connector entry). This shifts every row from the old
# row 16 onward down by two, which matches the target layout where the
# status-LED section's old rows 16-27 become rows 18-29.
$ws.Rows("16:17").Insert()

# Row 17 (the second of the two freshly inserted blank rows) becomes the
# new BJT transistor line item for the status LED driver circuit.
$ws.Range("C17").Value = "863-MMBT3904LT1G"
$ws.Range("A17").Value = "BJT transistor"
$ws.Range("B17").Value = "General purpose BJT transistor"
$ws.Range("D17").Value = 0.096
$ws.Range("E17").Value = 10
$ws.Range("F17").Formula = "=D17*E17"

# Match the number formatting used throughout column D/F for line items
# (reuse the existing currency-style cell format rather than defining a
# new one).
$ws.Range("D19").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("F19").Copy()
$ws.Range("F17").PasteSpecial(-4122)

[void]$ws.Range("F18").Select()
